$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104, pushing the current row 104 (and below) down to 105.
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new weekly record.
$ws.Range("A104").Value = 6
$ws.Range("B104").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C104").Value = "Metropolitana"
$ws.Range("D104").Value2 = 44448
$ws.Range("E104").Value = 13
$ws.Range("F104").Value = 100112029
$ws.Range("G104").Value = "Orégano"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 32
$ws.Range("K104").Value = 8500
$ws.Range("L104").Value = 9000
$ws.Range("M104").Value = 8734
$ws.Range("N104").Value = '$/docena de atados'
$ws.Range("O104").Value = "Región Metropolitana"
$ws.Range("P104").Value = 2911
$ws.Range("Q104").Value = 3
$ws.Range("R104").Value = "Hortaliza"

$ws.Range("D104").NumberFormat = "YYYY-MM-DD HH:MM:SS"
